$wb = $excel.ActiveWorkbook

# --- Sheet: Memo_Verification_details ---
$ws1 = $wb.Worksheets.Item("Memo_Verification_details")
$ws1.Range("B2").Value = "TESTINV90008"
$ws1.Range("K2").Value = "V0"
$ws1.Range("O2").Value = "test"
$ws1.Range("N2").Select()

# --- Sheet: Memo_invoice_Details ---
$ws2 = $wb.Worksheets.Item("Memo_invoice_Details")
$ws2.Range("B2").Value = "TESTINV90008"
$ws2.Range("C2").Value = "2024-03-13"
$ws2.Range("I2").Value = "0"
$ws2.Range("J2").Value = "0"
$ws2.Range("K2").Value = "0"
$ws2.Range("M2").Value = "test"
$ws2.Range("N2").Value = "test"
$ws2.Range("O2").Value = "1"

# Restore the horizontal alignment on cells that were right-aligned originally
# (I2 and O2 alignment swapped from left to right as part of this edit)
$ws2.Range("I2").HorizontalAlignment = -4152
$ws2.Range("O2").HorizontalAlignment = -4152
